# Update Spark memory diagram (single-slide deck).
#
# NOTE on numeric literals: PowerPoint's COM object model exposes shape
# geometry in points (Left/Top/Width/Height), while the underlying OOXML
# stores EMU (1 pt = 12700 EMU) as integers. The host's point -> EMU
# conversion round-trips through a 32-bit float, so naive
# `emu / 12700.0` literals can land 1 EMU off after truncation. The
# literals below were solved so that, after that float32 round trip,
# they reproduce the exact target EMU values from the diff.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Group 2 ("Off-Heap Unified Memory Pool" group): reposition ---
$group2 = $s.Shapes.Item(3)
$group2.Left = 609.0399780598426   # -> 7734807 EMU
$group2.Top  = 239.6500397         # -> 3043555 EMU

# Text tweak inside Group 2: "TextBox 104" (9th item in the group)
$tb104 = $group2.GroupItems.Item(9)
$tb104.TextFrame.TextRange.Text = "Optional Off-Heap Memory"
# Re-editing the text re-triggers autofit layout on this spAutoFit box;
# the diff does not touch its height, so pin it back explicitly.
$tb104.Height = 29.091378182677165  # -> 369460 EMU

# --- "Down Arrow 36": reposition + resize ---
$downArrow = $s.Shapes.Item(8)
$downArrow.Left   = 136.93468476929135  # -> 1739070 EMU
$downArrow.Top    = 268.36492918976376  # -> 3408234 EMU
$downArrow.Width  = 27.624054948031496  # -> 350825 EMU
$downArrow.Height = 42.93626027244095   # -> 545290 EMU

# --- Group 18: reposition ---
$group18 = $s.Shapes.Item(9)
$group18.Left = 45.97318840629921   # -> 583859 EMU
$group18.Top  = 316.98594671181104  # -> 4025721 EMU

# --- "Straight Connector 6": lengthen (position unchanged) ---
$conn6 = $s.Shapes.Item(10)
$conn6.Height = 399.1043701086614   # -> 5068625 EMU

# --- "Straight Arrow Connector 11": resize (position unchanged) ---
$conn11 = $s.Shapes.Item(13)
$conn11.Width  = 124.83846669685039  # -> 1585448 EMU
$conn11.Height = 8.453818827559054   # -> 107363 EMU

# --- "TextBox 52": reposition + add a second paragraph ---
$tb52 = $s.Shapes.Item(14)
[void]$tb52.TextFrame.TextRange.InsertAfter("`rOn-heap user memory size= (spark.executor.memory - 300 MB) * (1 - spark.memory.fraction)")
$tb52.Left = 609.302185084252    # -> 7738137 EMU
$tb52.Top  = 423.27003479999996  # -> 5375529 EMU

# --- "Straight Arrow Connector 54": reposition + resize ---
$conn54 = $s.Shapes.Item(15)
$conn54.Left   = 288.3540497480315   # -> 3662096 EMU
$conn54.Top    = 388.9551544102362   # -> 4939730 EMU
$conn54.Width  = 315.3699646598425   # -> 4005198 EMU
$conn54.Height = 44.74956705905512   # -> 568319 EMU
